$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency Price (D) and Volume(1h) (E) columns with the latest scraped values.
# A leading apostrophe forces Excel to store the value as text (matching the source data,
# which holds these figures as strings, e.g. "3.048.92" or "  -1.27%  ").
# The style is reset to "Normal" afterward so no extra number-format/quote-prefix style
# gets attached to the cell (keeping cell formatting identical to the original).

$ws.Range("D2").Value = "'62.892.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.19%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'3.042.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -1.54%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.18%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'581.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.89%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'150.39"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.09%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  -0.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.98%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'3.039.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.46%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -3.67%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'5.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -1.12%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.446"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.84%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.0000232"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.84%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'35.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -4.62%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +1.77%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.543.75"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -1.46%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'7.09"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.52%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'62.817.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.19%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.047.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.34%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'478.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.36%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -3.22%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.701"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -2.44%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'7.46"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -1.43%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.37"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.91%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'81.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.36%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'12.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.17%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'10.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +4.58%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'0.998"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -0.10%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'7.32"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.94%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.15%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'2.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.10%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -0.63%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'27.55"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +0.93%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  -3.89%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'1.05"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.10%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.0₃0804"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -5.55%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'5.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -4.14%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'2.17"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -2.51%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'3.11"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -8.49%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'50.27"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -0.89%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'9.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -2.55%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'423.64"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -4.87%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +3.14%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.283"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.68%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.832.46"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +1.01%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.0358"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.33%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'37.58"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -6.20%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'126.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -3.72%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +0.03%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'24.85"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -3.94%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -1.59%  "
$ws.Range("E51").Style = "Normal"
